$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 287892
$ws.Range("B2").Value = "CALL CENTER IN"
$ws.Range("C2").Value = "CREAR"
$ws.Range("F2").Value = "GESTOR"
$ws.Range("G2").Value = "ALLUS CALL INBOUND"
$ws.Range("H2").Value = "ALLUS INBOUND"
$ws.Range("I2").Value = "'+ SIMPLE"

# Row 3
$ws.Range("A3").Value = 287892
$ws.Range("B3").Value = "CALL CENTER IN"
$ws.Range("C3").Value = "CREAR"
$ws.Range("F3").Value = "GESTOR"
$ws.Range("G3").Value = "ALLUS CALL INBOUND"
$ws.Range("H3").Value = "ALLUS INBOUND"
$ws.Range("I3").Value = "'+ SIMPLE"

# Row 4
$ws.Range("A4").Value = 287892
$ws.Range("B4").Value = "CALL CENTER IN"
$ws.Range("C4").Value = "CREAR"
$ws.Range("E4").Value = "Sandra Vidal"
$ws.Range("G4").Value = "ALLUS CALL INBOUND"
$ws.Range("H4").Value = "ALLUS INBOUND II"
$ws.Range("I4").Value = "'+ SIMPLE"

# Match resulting selection state seen in target file
$ws.Range("F4").Select()
